# "Se crean nuevos metodos para automatizar epos"
#
# Refresh the MSISDN/SERIAL/MSI test fixtures on rows 12-13 of the data
# sheet, drop the old 4th sample row (row 14), and leave the cursor on
# the now-empty row beneath the table (row 14/15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: new SERIAL (col B) / PLU (col C) sample values.
# (Column B must stay text - these are 19/20-digit IMEI-style numbers -
#  and column C's 10-digit PLU also needs to round-trip as text, so we
#  assign them as strings rather than numbers.)
$ws.Range("B12").Value = "8957732111198172293"

# Row 13: new SERIAL / PLU / MSI sample values.
$ws.Range("B13").Value = "8957732111198172292"

$ws.Range("C12").Value = "3016876876"
$ws.Range("C13").Value = "3016877591"

$ws.Range("D13").Value = "732111198172292"

# The old row 14 sample (duplicate of row 12's MSI) is removed outright -
# not just cleared - so everything below it shifts up.
$ws.Rows.Item(14).Delete()

# A blank spacer row is left below the table with an explicit row height.
$ws.Rows.Item(15).RowHeight = 15

# Leave the selection / scroll position on the first empty row.
$ws.Range("A14").Select()
